$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column AK (column 37) header and values, mirroring existing format of column AJ (36)
$ws.Cells.Item(1, 37).Value = "28-jul"

$values = @{
    2  = 12
    3  = 16
    4  = 9
    5  = 9
    6  = 13
    7  = 14
    8  = 12
    9  = 13
    10 = 18
    11 = 15
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 37).Value = $values[$row]
}

# Copy formatting from column AJ (36) into new column AK (37) so styles match
$ws.Range($ws.Cells.Item(1, 36), $ws.Cells.Item(11, 36)).Copy() | Out-Null
$ws.Range($ws.Cells.Item(1, 37), $ws.Cells.Item(11, 37)).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Re-apply values after paste-special (paste only formats, values already set above) and update selection
$ws.Range("AK12").Select() | Out-Null

$wb.Save()
